# Loan RBI, Variable Instalments
#
# On the "Repayment schedule" sheet, a new (blank) column is inserted
# immediately before the "Late" column (column N). This pushes the
# existing "Late" column to O, the blank spacer/"heading" column to P,
# and "Outstanding" to Q -- matching the sheet growing from A:P to A:Q.
#
# It also records the workbook being left with the "Repayment schedule"
# tab active and cell M12 selected (previously the "Transactions" tab was
# the active one, with M7 selected on "Repayment schedule").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert the new blank column at position N (14th column), shifting
# "Late" / the spacer column / "Outstanding" one column to the right.
$ws.Columns.Item(14).Insert() | Out-Null

# Give the freshly inserted column an explicit (non best-fit) width,
# matching the neighbouring "In Advance" / "Due" style columns.
$ws.Columns.Item(14).ColumnWidth = 9.83

# Make "Repayment schedule" the active sheet/tab and select M12 (the new
# selection recorded in the workbook after the edit).
$ws.Activate() | Out-Null
$ws.Range("M12").Select() | Out-Null
